$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row with new columns I (I0) and J (IF), copying the
# existing header formatting (bold font + border + center/top alignment)
# from H1 so the new header cells match the style of the others.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the I0/IF data values for rows 2-76.
$data = @{
  2 = @(9, 9)
  3 = @(9, 9)
  4 = @(8, 8)
  5 = @(9, 9)
  6 = @(9, 9)
  7 = @(9, 9)
  8 = @(9, 9)
  9 = @(9, 9)
  10 = @(7, 7)
  11 = @(8, 8)
  12 = @(8, 8)
  13 = @(9, 9)
  14 = @(9, 9)
  15 = @(9, 9)
  16 = @(9, 9)
  17 = @(8, 9)
  18 = @(9, 9)
  19 = @(9, 9)
  20 = @(8, 9)
  21 = @(10, 10)
  22 = @(8, 9)
  23 = @(9, 9)
  24 = @(10, 10)
  25 = @(9, 9)
  26 = @(8, 9)
  27 = @(9, 9)
  28 = @(9, 9)
  29 = @(9, 9)
  30 = @(9, 9)
  31 = @(8, 9)
  32 = @(9, 9)
  33 = @(9, 9)
  34 = @(9, 9)
  35 = @(8, 9)
  36 = @(8, 9)
  37 = @(9, 9)
  38 = @(9, 9)
  39 = @(7, 7)
  40 = @(9, 9)
  41 = @(9, 9)
  42 = @(8, 8)
  43 = @(9, 9)
  44 = @(9, 9)
  45 = @(9, 9)
  46 = @(9, 9)
  47 = @(9, 9)
  48 = @(9, 9)
  49 = @(9, 9)
  50 = @(9, 9)
  51 = @(9, 9)
  52 = @(7, 7)
  53 = @(9, 9)
  54 = @(9, 9)
  55 = @(9, 9)
  56 = @(9, 9)
  57 = @(7, 7)
  58 = @(8, 9)
  59 = @(9, 9)
  60 = @(9, 9)
  61 = @(9, 10)
  62 = @(8, 8)
  63 = @(7, 7)
  64 = @(9, 9)
  65 = @(8, 9)
  66 = @(9, 9)
  67 = @(7, 7)
  68 = @(9, 9)
  69 = @(9, 9)
  70 = @(8, 8)
  71 = @(8, 9)
  72 = @(9, 9)
  73 = @(6, 6)
  74 = @(5, 5)
  75 = @(9, 9)
  76 = @(5, 5)
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Cells.Item($row, 9).Value = $vals[0]
  $ws.Cells.Item($row, 10).Value = $vals[1]
}
